# Adding negative tests, demo file and updating readme
# Edits the "ShareSkill" test-data sheet: drops the old blank/invalid rows
# (5-7), rewrites rows 2-4 with updated descriptions/negative-test data,
# switches the Start/End date columns to text-formatted dates, and
# centers the Skill Trade / Skill-Exchange columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ShareSkill")

# ---------------------------------------------------------------------
# Drop the old rows 5-7 (missing-title / unselected-subcategory / invalid
# file type negative tests that got replaced by the new negative tests
# baked directly into rows 2-4 below).
# ---------------------------------------------------------------------
$ws.Rows("5:7").Delete()

# ---------------------------------------------------------------------
# Row 2 - Software Testing: refresh the description text and store the
# start/end dates as literal (text) dates instead of Excel date serials.
# ---------------------------------------------------------------------
$ws.Range("B2").Value = "Test Automation using Selenium with C Sharp along with other testing frameworks and tools."

$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "'8/10/22022"
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "30/10/2022"

# ---------------------------------------------------------------------
# Row 3 - negative test: malicious/XSS title + description, work sample
# renamed, dates updated.
# ---------------------------------------------------------------------
$ws.Range("A3").Value = 'hack test <script>alert("YOU GOT HACKED by Javascript injection!");</script>'
$ws.Range("B3").Value = "Injecting Javascript or Malicious title along with special or invalid special characters"
$ws.Range("N3").Value = "wsample.png"

$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "'16/10/2022"
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "7/11/2022"

# ---------------------------------------------------------------------
# Row 4 - Copy Writer: negative test description (invalid upload),
# service type flipped to one-off, new work sample file, new dates/times.
# ---------------------------------------------------------------------
$ws.Range("B4").Value = "Upload invalid file type and size"
$ws.Range("F4").Value = "One-off service"
$ws.Range("N4").Value = "InvalidTypeAndSize.mp4"

$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "'8/10/2022"
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "16/10/2022"

$ws.Range("J4").Value = 0.3125
$ws.Range("K4").Value = 0.6875

# ---------------------------------------------------------------------
# Formatting touch-ups: header date cells share the header fill, Skill
# Trade / Skill-Exchange values are centered, row heights follow the new
# (longer) wrapped text.
# ---------------------------------------------------------------------
$ws.Range("H1:I1").NumberFormat = "@"

$ws.Range("L2:M4").HorizontalAlignment = -4108

$ws.Rows("2:2").RowHeight = 57.6
$ws.Rows("3:3").RowHeight = 86.4
$ws.Rows("4:4").RowHeight = 28.8

$ws.Range("I8").Select()
